$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"

$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("H2").Value = "2016-03-31 07:56:05"
$ws2.Range("H3").Value = "2016-03-31 07:56:05"

$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("H2").Value = "2016-03-31 07:56:22"
$ws3.Range("H3").Value = "2016-03-31 07:56:22"
